$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect new TPM-based calculations
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.05308866666666667
$ws.Range("H2").Value = 0.159266

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.05308866666666667
$ws.Range("N2").Value = 0.159266

$ws.Range("Q2").Value = 0.002818406528444445
$ws.Range("R2").Value = 0.025365658756
